# Update the "dSF" (column F) values for the rows that were repulled / recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -2
    4  = -3
    7  = -2
    9  = -1
    11 = 4
    13 = -6
    17 = -3
    19 = 2
    23 = 0
    25 = 1
    29 = 0
    30 = 1
    32 = -2
    34 = 6
    39 = -2
    40 = 1
    45 = -1
    51 = -2
    52 = -1
    53 = -1
    55 = -1
    56 = -1
    58 = -4
    64 = -2
    67 = -1
    69 = -2
    70 = 1
    71 = 1
    74 = -6
    75 = -6
    78 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
